# Auto-generated script applying scheduled-runner market data refresh
# to the Leve profit sheets (columns H-N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 9044
$ws.Range("I34").Value = 9044
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 9044
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -8841
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 9044
$ws.Range("I36").Value = 9044
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 9044
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -8329
$ws.Range("N36").ClearContents()
$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 3000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2752
$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 3000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2142
$ws.Range("H132").Value = 2156.625
$ws.Range("I132").Value = 2043
$ws.Range("K132").Value = 6129
$ws.Range("M132").Value = -3599
$ws.Range("H137").Value = 1295.3334
$ws.Range("J137").Value = 1178.6666
$ws.Range("L137").Value = 3535.9998
$ws.Range("N137").Value = -8635.9998
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 918.1667
$ws.Range("I4").Value = 927
$ws.Range("K4").Value = 927
$ws.Range("M4").Value = -811
$ws.Range("H61").Value = 4197.5835
$ws.Range("I61").Value = 3264.7778
$ws.Range("K61").Value = 3264.7778
$ws.Range("M61").Value = -3052.7778
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3669.3635
$ws.Range("I132").Value = 3590
$ws.Range("K132").Value = 10770
$ws.Range("M132").Value = -8240
$ws.Range("H136").Value = 4197.5835
$ws.Range("I136").Value = 3264.7778
$ws.Range("K136").Value = 9794.3334
$ws.Range("M136").Value = -7244.3334
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 499.33334
$ws.Range("I22").Value = 499.33334
$ws.Range("K22").Value = 499.33334
$ws.Range("M22").Value = -326.33334
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2100.9048
$ws.Range("I31").Value = 1861.8125
$ws.Range("J31").Value = 2866
$ws.Range("K31").Value = 1861.8125
$ws.Range("L31").Value = 2866
$ws.Range("M31").Value = -1566.8125
$ws.Range("N31").Value = -3456
$ws.Range("H34").Value = 2100.9048
$ws.Range("I34").Value = 1861.8125
$ws.Range("J34").Value = 2866
$ws.Range("K34").Value = 1861.8125
$ws.Range("L34").Value = 2866
$ws.Range("M34").Value = -1659.8125
$ws.Range("N34").Value = -3270
$ws.Range("H134").Value = 5850
$ws.Range("I134").Value = 5850
$ws.Range("K134").Value = 17550
$ws.Range("M134").Value = -15015
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 299
$ws.Range("I35").Value = 299
$ws.Range("K35").Value = 897
$ws.Range("M35").Value = -609
$ws.Range("H68").Value = 2550.975
$ws.Range("I68").Value = 918.0909
$ws.Range("J68").Value = 3170.3447
$ws.Range("K68").Value = 2754.2727
$ws.Range("L68").Value = 9511.034100000001
$ws.Range("M68").Value = -1943.2727
$ws.Range("N68").Value = -11133.0341
$ws.Range("H71").Value = 2550.975
$ws.Range("I71").Value = 918.0909
$ws.Range("J71").Value = 3170.3447
$ws.Range("K71").Value = 8262.8181
$ws.Range("L71").Value = 28533.1023
$ws.Range("M71").Value = -4206.8181
$ws.Range("N71").Value = -36645.1023
$ws.Range("H76").Value = 9017.5
$ws.Range("I76").Value = 9017.5
$ws.Range("K76").Value = 27052.5
$ws.Range("M76").Value = -26669.5
$ws.Range("H79").Value = 9017.5
$ws.Range("I79").Value = 9017.5
$ws.Range("K79").Value = 27052.5
$ws.Range("M79").Value = -25726.5
$ws.Range("H117").Value = 3229.5
$ws.Range("J117").Value = 3715.6
$ws.Range("L117").Value = 11146.8
$ws.Range("N117").Value = -18030.8
$ws.Range("H121").Value = 18136
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 90000
$ws.Range("N121").Value = -92620
$ws.Range("H122").Value = 2890.5
$ws.Range("I122").Value = 2248.5
$ws.Range("J122").Value = 3211.5
$ws.Range("K122").Value = 20236.5
$ws.Range("L122").Value = 28903.5
$ws.Range("M122").Value = -17786.5
$ws.Range("N122").Value = -33803.5
$ws.Range("H125").Value = 2400
$ws.Range("I125").Value = 2400
$ws.Range("K125").Value = 7200
$ws.Range("M125").Value = -2280
$ws.Range("H131").Value = 16754.465
$ws.Range("I131").Value = 102091.55
$ws.Range("K131").Value = 306274.65
$ws.Range("M131").Value = -301234.65
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2045.5
$ws.Range("I122").Value = 2591
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 7773
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -5323
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 6773.5
$ws.Range("I132").Value = 4634
$ws.Range("K132").Value = 13902
$ws.Range("M132").Value = -11372
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1853
$ws.Range("I22").Value = 1763.3636
$ws.Range("J22").Value = 2099.5
$ws.Range("K22").Value = 1763.3636
$ws.Range("L22").Value = 2099.5
$ws.Range("M22").Value = -1468.3636
$ws.Range("N22").Value = -2689.5
$ws.Range("H27").Value = 1853
$ws.Range("I27").Value = 1763.3636
$ws.Range("J27").Value = 2099.5
$ws.Range("K27").Value = 1763.3636
$ws.Range("L27").Value = 2099.5
$ws.Range("M27").Value = -1656.3636
$ws.Range("N27").Value = -2313.5
$ws.Range("H46").Value = 2701.125
$ws.Range("I46").Value = 1698.4
$ws.Range("J46").Value = 3156.9092
$ws.Range("K46").Value = 1698.4
$ws.Range("L46").Value = 3156.9092
$ws.Range("M46").Value = -1510.4
$ws.Range("N46").Value = -3532.9092
$ws.Range("H55").Value = 251.25
$ws.Range("I55").Value = 251.33333
$ws.Range("K55").Value = 251.33333
$ws.Range("M55").Value = -78.33332999999999
$ws.Range("H68").Value = 2997.7273
$ws.Range("I68").Value = 2947.5
$ws.Range("J68").Value = 3500
$ws.Range("K68").Value = 2947.5
$ws.Range("L68").Value = 3500
$ws.Range("M68").Value = -2198.5
$ws.Range("N68").Value = -4998
$ws.Range("H71").Value = 2997.7273
$ws.Range("I71").Value = 2947.5
$ws.Range("J71").Value = 3500
$ws.Range("K71").Value = 14737.5
$ws.Range("L71").Value = 17500
$ws.Range("M71").Value = -10993.5
$ws.Range("N71").Value = -24988
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 36014.668
$ws.Range("I39").Value = 26522
$ws.Range("K39").Value = 26522
$ws.Range("M39").Value = -26109
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H132").Value = 4517.25
$ws.Range("J132").Value = 8488.6
$ws.Range("L132").Value = 25465.8
$ws.Range("N132").Value = -30525.8
